$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pilot_points")

# Update mean_kh values (column E) for several pilot points
$ws.Range("E7").Value = 2
$ws.Range("E8").Value = 20
$ws.Range("E9").Value = 5
$ws.Range("E10").Value = 5
$ws.Range("E24").Value = 0.5
$ws.Range("E26").Value = 10
$ws.Range("E27").Value = 2
$ws.Range("E42").Value = 1

# Update the view: activate the sheet and move the selection/scroll position
$ws.Activate()
$ws.Range("H6").Select()
